$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - increment "想去人数" (want-to-go count) in column F
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 190
$ws1.Range("F4").Value = 789
$ws1.Range("F5").Value = 70
$ws1.Range("F6").Value = 9

# Sheet "全部类型" (all types) - same events repeated, increment the same counts
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 190
$ws4.Range("F5").Value = 789
$ws4.Range("F6").Value = 70
$ws4.Range("F7").Value = 9
